$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 147
$ws1.Range("F6").Value = 3156
$ws1.Range("F8").Value = 595
$ws1.Range("F10").Value = 645
$ws1.Range("F12").Value = 542
$ws1.Range("F13").Value = 404
$ws1.Range("F16").Value = 1375
$ws1.Range("F18").Value = 1641
$ws1.Range("F21").Value = 616
$ws1.Range("F27").Value = 110
$ws1.Range("F30").Value = 40
$ws1.Range("F32").Value = 3978
$ws1.Range("F36").Value = 1437
$ws1.Range("F38").Value = 1884

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 56

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 147
$ws4.Range("F6").Value = 3156
$ws4.Range("F8").Value = 595
$ws4.Range("F10").Value = 645
$ws4.Range("F12").Value = 542
$ws4.Range("F14").Value = 404
$ws4.Range("F17").Value = 1375
$ws4.Range("F19").Value = 1641
$ws4.Range("F22").Value = 616
$ws4.Range("F28").Value = 110
$ws4.Range("F31").Value = 40
$ws4.Range("F33").Value = 3978
$ws4.Range("F34").Value = 56
$ws4.Range("F39").Value = 1437
$ws4.Range("F41").Value = 1884
